# Assertion added for HealtRisk Assessment
#
# - Rename Sheet1 -> UserData
# - Add a new sheet "HealtAssessmentData" after UserData, becomes the active tab
# - New sheet gets headers Completed_Status / InCompleted_Status and a status row
# - UserData!C2 keeps its value ("Anas") but the view selection / tabSelected move

$wb = $excel.ActiveWorkbook

$ws1 = $wb.ActiveSheet
$ws1.Name = "UserData"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "HealtAssessmentData"

$ws2.Range("A1").Value = "Completed_Status"
$ws2.Range("B1").Value = "InCompleted_Status"
$ws2.Range("A2").Value = "completed on 04.01.2016"
$ws2.Range("A2").Style = "Hyperlink"

$ws2.Columns.Item(1).ColumnWidth = 22.5
$ws2.Columns.Item(2).ColumnWidth = 18.5

$ws2.Range("B8").Select() | Out-Null

$ws1.Range("A15").Select() | Out-Null

$ws2.Activate() | Out-Null
